$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns L1:N1
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copy header style from existing header cell (e.g. A1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update E and F columns (percent columns) - multiply existing fraction by 100
# since number format stays 0.00% but underlying values change to already-scaled percentages
$ws.Range("E2").Value = 85.0187265917603
$ws.Range("F2").Value = 60.79295154185021

$ws.Range("E3").Value = 14.9812734082397
$ws.Range("F3").Value = 70

$ws.Range("E4").Value = 80.17711171662125
$ws.Range("F4").Value = 93.28802039082413

$ws.Range("E5").Value = 19.82288828337875
$ws.Range("F5").Value = 97.9381443298969

$ws.Range("E6").Value = 96.49122807017544
$ws.Range("F6").Value = 21.36363636363636

$ws.Range("E7").Value = 3.508771929824561
$ws.Range("F7").Value = 45.83333333333333

# New columns L, M, N data
$ws.Range("L2").Value = 92.42558549437925
$ws.Range("M2").Value = 211586
$ws.Range("N2").Value = 306.6463768115942

$ws.Range("L3").Value = 86.73427165141923
$ws.Range("M3").Value = 51967
$ws.Range("N3").Value = 371.1928571428571

$ws.Range("L4").Value = 88.21289274334562
$ws.Range("M4").Value = 150490
$ws.Range("N4").Value = 137.0582877959927

$ws.Range("L5").Value = 95.69895606541161
$ws.Range("M5").Value = 53156
$ws.Range("N5").Value = 186.5122807017544

$ws.Range("L6").Value = 19.42841748106926
$ws.Range("M6").Value = 2114
$ws.Range("N6").Value = 14.99290780141844

$ws.Range("L7").Value = 22.50301276387648
$ws.Range("M7").Value = 94
$ws.Range("N7").Value = 8.545454545454545

$wb.Save()
